$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rebuild the device/feature table to match the new feature model.
# A new "Devices" column is inserted at the front (column A); the remaining
# columns (Input, Language, Output, Contacts, Alerts, Disabilities) keep
# their relative order and a new "Patterns" column is appended at the end.
# ---------------------------------------------------------------------------

# Header row
$ws.Range("A2").Value = "Devices"
$ws.Range("B2").Value = "Input"
$ws.Range("C2").Value = "Language"
$ws.Range("D2").Value = "Output"
$ws.Range("E2").Value = "Contacts"
$ws.Range("F2").Value = "Alerts"
$ws.Range("G2").Value = "Disabilities"
$ws.Range("H2").Value = "Patterns"

# Row 3 - Portable Alarm
$ws.Range("A3").Value = "Portable Alarm"
$ws.Range("B3").Value = "Mouse, Keyboard"
$ws.Range("C3").Value = "PT"
$ws.Range("D3").Value = "Screen, Light"
$ws.Range("E3").Value = "Management"
$ws.Range("F3").Value = "SMS"
$ws.Range("G3").Value = "Auditory"
$ws.Range("H3").Value = "Reminders"

# Row 4 - Smart Watch
$ws.Range("A4").Value = "Smart Watch"
$ws.Range("B4").Value = "Mouse, Keyboard, Microphone"
$ws.Range("C4").Value = "EN"
$ws.Range("D4").Value = "Screen, Speakers"
$ws.Range("E4").Value = "Management"
$ws.Range("F4").Value = "SMS"
$ws.Range("G4").Value = "Visual"
$ws.Range("H4").Value = "Reminders"

# Row 5 - Panic Vest
$ws.Range("A5").Value = "Panic Vest"
$ws.Range("B5").Value = "Mouse, Keyboard"
$ws.Range("C5").Value = "EN"
$ws.Range("D5").Value = "Screen, Light"
$ws.Range("E5").Value = "Management, Emergency Contact"
$ws.Range("F5").Value = "SMS"
$ws.Range("G5").Value = "Auditory"
$ws.Range("H5").Value = "Reminders"

# Row 6 - Movement Sensor
$ws.Range("A6").Value = "Movement Sensor"
$ws.Range("B6").Value = "Mouse, Keyboard"
$ws.Range("C6").Value = "PT"
$ws.Range("D6").Value = "Screen, Speakers"
$ws.Range("E6").Value = "Management"
$ws.Range("F6").Value = "SMS"
$ws.Range("G6").Value = "Visual"
$ws.Range("H6").Value = "Reminders, Movement Detection"

# Column H is no longer used on the blank trailer rows (7 and 8); the table
# now only spans A:G there, so drop the cell entirely (not just its value).
$ws.Range("H7").Delete(-4159)
$ws.Range("H8").Delete(-4159)

# Make sure the blank trailer rows have the same formatted-but-empty cell in
# the new column A as the rest of the table (they were shifted from "none").
$ws.Range("A7").Value = ""
$ws.Range("A8").Value = ""

# ---------------------------------------------------------------------------
# Formatting: every data cell in A2:H8 (except B3:B6, see below) is centered
# with the plain default font. Column F used to show a decorative Webdings
# glyph ("a") as a stand-in icon; it is now plain centered text, so drop the
# special font there. B3:B6 already carry the desired centered styling and
# are left untouched.
# ---------------------------------------------------------------------------
$ws.Range("A2:H2").HorizontalAlignment = -4108
$ws.Range("A3:A6").HorizontalAlignment = -4108
$ws.Range("C3:H6").HorizontalAlignment = -4108
$ws.Range("A7:G8").HorizontalAlignment = -4108

$ws.Range("F3:F6").Value = "SMS"
$ws.Range("F3:F6").Style = "Normal"
$ws.Range("F3:F6").HorizontalAlignment = -4108

# Row 3-6 no longer carry an explicit 15.75 row height / 0.3 dyDescent - they
# now share the sheet's default row height like row 2.
$ws.Rows("3:6").RowHeight = 15

# ---------------------------------------------------------------------------
# Column widths - auto fit to the new (generally longer) content.
# ---------------------------------------------------------------------------
$ws.Columns("A:H").AutoFit()

# ---------------------------------------------------------------------------
# View state: zoom in a bit and leave the selection on the last data column.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 140
$ws.Range("H7").Select()
